$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.459.59'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '2.985.09'
$ws.Range('E3').Value = '  +2.60%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '382.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.81%  '
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.594'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.05%  '
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.139'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('E12').Value = '  +1.61%  '
$ws.Range('D13').Value = '3.453.78'
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.55'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.86%  '
$ws.Range('D16').Value = '2.995.54'
$ws.Range('E16').Value = '  +3.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.972'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.59%  '
$ws.Range('D18').Value = '51.419.62'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('E19').Value = '  +3.68%  '
$ws.Range('E20').Value = '  +3.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('E22').Value = '  +2.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '262.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.37'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +18.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +24.38%  '
$ws.Range('E28').Value = '  +15.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.171'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.44%  '
$ws.Range('E30').Value = '  +1.84%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '50.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0452'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.70%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.04'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.07'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('E42').Value = '  +3.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '122.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('E45').Value = '  +18.26%  '
$ws.Range('E46').Value = '  -2.42%  '
$ws.Range('E47').Value = '  +2.92%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.032.04'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.27'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0334'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.90%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.29'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.14%  '
